$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 0.04132999999999999
$ws.Range("E2").Value = -0.09745000000000001
$ws.Range("F2").Value = 0.09960000000000001
$ws.Range("G2").Value = 0.07598466291014208
$ws.Range("H2").Value = 0.07598466291014208
$ws.Range("I2").Value = 0.04995108215559357
$ws.Range("J2").Value = 0.0406194526638052
$ws.Range("K2").Value = 2546.7
$ws.Range("L2").Value = 0.02687790169360226
$ws.Range("M2").Value = 2142.15
$ws.Range("N2").Value = 0.03521442802753831
$ws.Range("O2").Value = 0.8411473671810578
$ws.Range("P2").Value = 2142.15
$ws.Range("Q2").Value = 0.03521442802753831
$ws.Range("R2").Value = 0.8411473671810578
$ws.Range("U2").Value = 8115.700000000001
$ws.Range("V2").Value = 0.1334125684677043
$ws.Range("W2").Value = 0.06849407011234906
$ws.Range("X2").Value = 0.07582208041715149
$ws.Range("Y2").Value = -0.007328010304802429
$ws.Range("Z2").Value = 2.197341873721609
$ws.Range("AA2").Value = 0.09302589993455763
$ws.Range("AB2").Value = 0.0662963701743231
$ws.Range("AC2").Value = 0.02672952976023453
$ws.Range("AD2").Value = 10231.9
$ws.Range("AE2").Value = 0
$ws.Range("AF2").Value = 10231.9
$ws.Range("AG2").Value = 2116.199999999999
$ws.Range("AH2").Value = 0.1439824945295405
$ws.Range("AI2").Value = 0.1742822200419699
$ws.Range("AJ2").Value = 0.03361833137933334
$ws.Range("AK2").Value = 0.04182783818346769
$ws.Range("AL2").Value = 672.7
$ws.Range("AM2").Value = 672.7
$ws.Range("AN2").Value = 1.984926670287887
$ws.Range("AO2").Value = 7.035677122045487
$ws.Range("AP2").Value = 0.4105299914642661
$ws.Range("AQ2").Value = 7.035677122045487
$ws.Range("D3").Value = 0.07519999999999999
$ws.Range("E3").Value = -0.0249
$ws.Range("F3").Value = 0.0762
$ws.Range("G3").Value = 0.08288043478260869
$ws.Range("H3").Value = 0.08288043478260869
$ws.Range("I3").Value = 0.06125701743908265
$ws.Range("J3").Value = 0.04824886597644338
$ws.Range("K3").Value = 1114.1
$ws.Range("L3").Value = 0.04158579192546583
$ws.Range("M3").Value = 531
$ws.Range("N3").Value = 0.02762488424601234
$ws.Range("O3").Value = 0.4766178978547707
$ws.Range("P3").Value = 531
$ws.Range("Q3").Value = 0.02762488424601234
$ws.Range("R3").Value = 0.4766178978547707
$ws.Range("U3").Value = 1578.6
$ws.Range("V3").Value = 0.08212550333475532
$ws.Range("W3").Value = 0.09522466388027043
$ws.Range("X3").Value = 0.07708596959357696
$ws.Range("Y3").Value = 0.01813869428669347
$ws.Range("Z3").Value = 2.07814451382694
$ws.Range("AA3").Value = 0.1002681161273171
$ws.Range("AB3").Value = 0.0659915570069101
$ws.Range("AC3").Value = 0.03427655912040702
$ws.Range("AD3").Value = 4036.5
$ws.Range("AE3").Value = 0
$ws.Range("AF3").Value = 4036.5
$ws.Range("AG3").Value = 2457.9
$ws.Range("AH3").Value = 0.1735509474037225
$ws.Range("AI3").Value = 0.2279439584827454
$ws.Range("AJ3").Value = 0.1133733400369931
$ws.Range("AK3").Value = 0.1523834913234592
$ws.Range("AL3").Value = 233.5
$ws.Range("AM3").Value = 233.5
$ws.Range("AN3").Value = 2.238520408163265
$ws.Range("AO3").Value = 7.028265524625267
$ws.Range("AP3").Value = 1.363076752440106
$ws.Range("AQ3").Value = 7.028265524625267
$ws.Range("B4").Value = "Münchener Rückversicherungs-Gesellschaft Aktiengesellschaft in München (XTRA:MUV2)"
$ws.Range("D4").Value = 0.00746
$ws.Range("E4").Value = -0.17
$ws.Range("F4").Value = 0.123
$ws.Range("G4").Value = 0.07326630400395524
$ws.Range("H4").Value = 0.07326630400395524
$ws.Range("I4").Value = 0.0454942076477002
$ws.Range("J4").Value = 0.0381570351249066
$ws.Range("K4").Value = 1432.6
$ws.Range("L4").Value = 0.02107995403198632
$ws.Range("M4").Value = 1611.15
$ws.Range("N4").Value = 0.03872044566424255
$ws.Range("O4").Value = 1.124633533435711
$ws.Range("P4").Value = 1611.15
$ws.Range("Q4").Value = 0.03872044566424255
$ws.Range("R4").Value = 1.124633533435711
$ws.Range("U4").Value = 6537.1
$ws.Range("V4").Value = 0.1571048166537691
$ws.Range("W4").Value = 0.0417634763444277
$ws.Range("X4").Value = 0.07455819124072605
$ws.Range("Y4").Value = -0.03279471489629834
$ws.Range("Z4").Value = 2.248174771991227
$ws.Range("AA4").Value = 0.08578368374179815
$ws.Range("AB4").Value = 0.06660118334173611
$ws.Range("AC4").Value = 0.01918250040006204
$ws.Range("AD4").Value = 6195.4
$ws.Range("AE4").Value = 0
$ws.Range("AF4").Value = 6195.4
$ws.Range("AG4").Value = -341.7000000000007
$ws.Range("AH4").Value = 0.1295967802665819
$ws.Range("AI4").Value = 0.1511054743234838
$ws.Range("AJ4").Value = -0.008280003198596512
$ws.Range("AK4").Value = -0.009914866205888008
$ws.Range("AL4").Value = 439.2
$ws.Range("AM4").Value = 439.2
$ws.Range("AN4").Value = 1.848490273302303
$ws.Range("AO4").Value = 7.039617486338798
$ws.Range("AP4").Value = -0.1019513068385251
$ws.Range("AQ4").Value = 7.039617486338798
